# Add a new "12-aug" column (BH) to the "Prix Spot" sheet, one column to
# the right of the existing last column "11-aug" (BG): a new header cell
# in row 1 plus 24 hourly price values in rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# BG is column 59 (the current last column, "11-aug"); BH (60) is new.
$prevCol = 59
$newCol = 60

# Header cell: copy the previous header's formatting (bold, border,
# centered) onto the new cell, then set its text.
$srcHeader = $ws.Cells.Item(1, $prevCol)
$dstHeader = $ws.Cells.Item(1, $newCol)
$srcHeader.Copy($dstHeader)
$dstHeader.Value = "12-aug"

# Hourly values for rows 2-25 ("12-aug" column).
$values = @{
    2  = 93
    3  = 89.67
    4  = 85.41
    5  = 79.06
    6  = 69.67
    7  = 66.34
    8  = 75.52
    9  = 102.81
    10 = 98.69
    11 = 92.5
    12 = 66.34
    13 = 40.94
    14 = 8.380000000000001
    15 = 3
    16 = 5.15
    17 = 35.05
    18 = 63.8
    19 = 85.81999999999999
    20 = 96.78
    21 = 143.64
    22 = 154.93
    23 = 140.81
    24 = 118.41
    25 = 100.37
}

foreach ($row in 2..25) {
    $ws.Cells.Item($row, $newCol).Value = $values[$row]
}
